$wb = $excel.ActiveWorkbook

# Update the status text from "Ready for handoff" to "In Translation"
# everywhere it appears (the Overview summary columns + the per-locale
# "Status" table column on each language sheet).
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        if ("Ready for handoff" -eq $cell.Value2) {
            $cell.Value = "In Translation"
        }
    }
}

# The columns that held the old (longer) status text shrink to fit the
# new, shorter text.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E1:F1").ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C1").ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C1").ColumnWidth = 12.5
